$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Cells.Item(19, 1).Value = 111670558
$ws.Cells.Item(19, 17).Value = 558133.6011735104
$ws.Cells.Item(19, 18).Value = 7067979.426396712

# Row 20
$ws.Cells.Item(20, 1).Value = 111671294
$ws.Cells.Item(20, 2).Value = 78578
$ws.Cells.Item(20, 4).Value = 'NT'
$ws.Cells.Item(20, 5).Value = 6458
$ws.Cells.Item(20, 6).Value = 'Lunglav'
$ws.Cells.Item(20, 7).Value = 'Lobaria pulmonaria'
$ws.Cells.Item(20, 8).Value = '(L.) Hoffm.'
$ws.Cells.Item(20, 17).Value = 558118.4535210516
$ws.Cells.Item(20, 18).Value = 7067742.103054954

# Row 21
$ws.Cells.Item(21, 1).Value = 111670497
$ws.Cells.Item(21, 17).Value = 558159.8619213518
$ws.Cells.Item(21, 18).Value = 7068022.886732788

# Row 22
$ws.Cells.Item(22, 1).Value = 111671226
$ws.Cells.Item(22, 2).Value = 78579
$ws.Cells.Item(22, 5).Value = 2081
$ws.Cells.Item(22, 6).Value = 'Skrovellav'
$ws.Cells.Item(22, 7).Value = 'Lobaria scrobiculata'
$ws.Cells.Item(22, 8).Value = '(Scop.) DC.'
$ws.Cells.Item(22, 17).Value = 558118.4535210516
$ws.Cells.Item(22, 18).Value = 7067742.103054954

# Row 23
$ws.Cells.Item(23, 1).Value = 111671190
$ws.Cells.Item(23, 2).Value = 78611
$ws.Cells.Item(23, 4).Value = 'LC'
$ws.Cells.Item(23, 5).Value = 6463
$ws.Cells.Item(23, 6).Value = 'Bårdlav'
$ws.Cells.Item(23, 7).Value = 'Nephroma parile'
$ws.Cells.Item(23, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(23, 17).Value = 558215.9329796816
$ws.Cells.Item(23, 18).Value = 7067869.292590594

# Row 24
$ws.Cells.Item(24, 1).Value = 111671179
$ws.Cells.Item(24, 17).Value = 558215.9656782644
$ws.Cells.Item(24, 18).Value = 7067867.520903144

# Row 25
$ws.Cells.Item(25, 1).Value = 111671188
$ws.Cells.Item(25, 2).Value = 78605
$ws.Cells.Item(25, 4).Value = 'LC'
$ws.Cells.Item(25, 5).Value = 6462
$ws.Cells.Item(25, 6).Value = 'Stuplav'
$ws.Cells.Item(25, 7).Value = 'Nephroma bellum'
$ws.Cells.Item(25, 8).Value = '(Spreng.) Tuck.'
$ws.Cells.Item(25, 17).Value = 558215.9329796816
$ws.Cells.Item(25, 18).Value = 7067869.292590594
$ws.Cells.Item(25, 12).Value = ""  # clears placeholder empty cell

# Row 26
$ws.Cells.Item(26, 1).Value = 111670477
$ws.Cells.Item(26, 2).Value = 96346
$ws.Cells.Item(26, 5).Value = 620
$ws.Cells.Item(26, 6).Value = 'Skogsfru'
$ws.Cells.Item(26, 7).Value = 'Epipogium aphyllum'
$ws.Cells.Item(26, 8).Value = 'Sw.'
$ws.Cells.Item(26, 17).Value = 558155.0815836267
$ws.Cells.Item(26, 18).Value = 7068017.481975557
$ws.Cells.Item(26, 12).Value = ""  # placeholder empty cell (was absent, diff expects empty inlineStr cell)

# Row 27
$ws.Cells.Item(27, 1).Value = 111670510
$ws.Cells.Item(27, 2).Value = 96346
$ws.Cells.Item(27, 4).Value = 'NT'
$ws.Cells.Item(27, 5).Value = 620
$ws.Cells.Item(27, 6).Value = 'Skogsfru'
$ws.Cells.Item(27, 7).Value = 'Epipogium aphyllum'
$ws.Cells.Item(27, 8).Value = 'Sw.'
$ws.Cells.Item(27, 17).Value = 558124.4538526792
$ws.Cells.Item(27, 18).Value = 7067994.321708324
$ws.Cells.Item(27, 12).Value = ""  # placeholder empty cell (was absent, diff expects empty inlineStr cell)

# Row 28
$ws.Cells.Item(28, 1).Value = 111671201
$ws.Cells.Item(28, 2).Value = 78579
$ws.Cells.Item(28, 5).Value = 2081
$ws.Cells.Item(28, 6).Value = 'Skrovellav'
$ws.Cells.Item(28, 7).Value = 'Lobaria scrobiculata'
$ws.Cells.Item(28, 8).Value = '(Scop.) DC.'
$ws.Cells.Item(28, 17).Value = 558250.1783714101
$ws.Cells.Item(28, 18).Value = 7067936.828089682
$ws.Cells.Item(28, 12).Value = ""  # clears placeholder empty cell

# Row 29
$ws.Cells.Item(29, 1).Value = 111670567
$ws.Cells.Item(29, 17).Value = 558129.9933989302
$ws.Cells.Item(29, 18).Value = 7067958.536170656

# Row 30
$ws.Cells.Item(30, 1).Value = 111671197
$ws.Cells.Item(30, 2).Value = 78578
$ws.Cells.Item(30, 5).Value = 6458
$ws.Cells.Item(30, 6).Value = 'Lunglav'
$ws.Cells.Item(30, 7).Value = 'Lobaria pulmonaria'
$ws.Cells.Item(30, 8).Value = '(L.) Hoffm.'
$ws.Cells.Item(30, 17).Value = 558250.1783714101
$ws.Cells.Item(30, 18).Value = 7067936.828089682
